# Apply "ultimo aggiornamento prima di nuovo setup fatto da Peter" update:
# fill in missing xG_away / goals_home / goals_away values for rows 2-9,
# and add the full xG_home / xG_away / goals_home / goals_away data for
# the newly-completed rows 10-15.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, [string]$value)
    # Force the cell to be written as text (matching the shared-string
    # cells already used throughout this sheet) instead of letting Excel
    # auto-convert numeric-looking strings into number cells.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Rows 2-9: xG_home (D) already present; update xG_away (E), goals_home (F),
# goals_away (G).
$updates = @{
    2 = @{ E = "0.95973";  F = "1"; G = "1" }
    3 = @{ E = "2.6497";   F = "1"; G = "4" }
    4 = @{ E = "1.5945";   F = "4"; G = "1" }
    5 = @{ E = "3.23668";  F = "3"; G = "4" }
    6 = @{ E = "2.08517";  F = "3"; G = "3" }
    7 = @{ E = "1.52493";  F = "0"; G = "2" }
    8 = @{ E = "0.0917939"; F = "0"; G = "0" }
    9 = @{ E = "0.51";     F = "0"; G = "2" }
}

foreach ($row in $updates.Keys | Sort-Object) {
    $vals = $updates[$row]
    Set-TextValue $ws.Range("E$row") $vals.E
    Set-TextValue $ws.Range("F$row") $vals.F
    Set-TextValue $ws.Range("G$row") $vals.G
}

# Rows 10-15: add brand-new xG_home (D), xG_away (E), goals_home (F),
# goals_away (G) values.
$newRows = @{
    10 = @{ D = "0.753035"; E = "1.49985";  F = "0"; G = "3" }
    11 = @{ D = "1.27993";  E = "0.209882"; F = "0"; G = "0" }
    12 = @{ D = "1.08927";  E = "2.09826";  F = "1"; G = "0" }
    13 = @{ D = "1.71835";  E = "1.7103";   F = "1"; G = "1" }
    14 = @{ D = "1.04718";  E = "1.41457";  F = "1"; G = "2" }
    15 = @{ D = "1.6455";   E = "2.16381";  F = "2"; G = "3" }
}

foreach ($row in $newRows.Keys | Sort-Object) {
    $vals = $newRows[$row]
    Set-TextValue $ws.Range("D$row") $vals.D
    Set-TextValue $ws.Range("E$row") $vals.E
    Set-TextValue $ws.Range("F$row") $vals.F
    Set-TextValue $ws.Range("G$row") $vals.G
}
